$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting
# existing "Late" / blank / "Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()

# Update the active selection to match the post-edit workbook state.
$ws.Range("S7").Select()
